# Reverted and updated VCEA and 1.5 scenarios and policy schedules; updated
# trans/RTMF to reflect primarily reduction in VMT with only 5% of VMT
# change going to passenger rail.
#
# This updates the "RTMF-passengers" worksheet of the Recipient
# Transportation Mode Fractions workbook:
#   - LDVs -> aircraft (C2):  0.33 -> 0.15
#   - LDVs -> rail     (E2):  0.33 -> 0.05
#   - Non-motorized/eliminated share (I2) is hard-coded to 0.8 instead of
#     being computed with the 1-SUM(B2:G2) formula.
# It also makes RTMF-passengers the active sheet/tab with the selection
# resting on cell E4, matching the saved view state of the workbook.

$wb = $excel.ActiveWorkbook

$wsPassengers = $wb.Worksheets.Item("RTMF-passengers")

# Update the LDV row (row 2): mode-shift fractions to aircraft and rail.
$wsPassengers.Range("C2").Value = 0.15
$wsPassengers.Range("E2").Value = 0.05

# The "non-motorized / eliminated" column no longer derives from the
# formula; it is now a plain, hard-coded value.
$wsPassengers.Range("I2").Formula = 0.8

# Make RTMF-passengers the active/selected sheet, with E4 selected,
# matching the workbook's saved view state.
$wsPassengers.Activate()
$wsPassengers.Range("E4").Select()
